# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-03 09:19:13
#
# Applies the attendance-data refresh to the "Session Analysis Results" sheet:
#  - Updated "Recorded By" lists / "Students" counts for several sessions
#    (re-ordered / newly added recorders, updated counts)
#  - Three sessions (rows 24, the A2/ANATOMY #1 session) flipped from
#    Pending -> Recorded, while four other sessions (rows 98, 117, 146, 172)
#    flipped from Pending -> Not Recorded, each carrying the matching
#    row-highlight color (green / pink) instead of the "Pending" yellow.
#  - Downstream Class Statistics / breakdown table numbers recomputed to
#    match the new raw attendance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (Year2/A1/ANATOMY #1): recorder list grew, student count updated
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("H2").Value = "163/216"

# ---------------------------------------------------------------------
# Class Statistics block (K/L columns)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 12
$ws.Range("L8").Value = 142
$ws.Range("L9").Value = "12.5%"
$ws.Range("L10").Value = "30.2%"

# ---------------------------------------------------------------------
# Per-group breakdown table (K:S columns), rows 15-22
# ---------------------------------------------------------------------
# Row 15 (A1)
$ws.Range("S15").Value = "36.3%"

# Row 16 (A2)
$ws.Range("O16").Value = 2
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "9.1%"
$ws.Range("S16").Value = "31.1%"

# Row 19 (B1)
$ws.Range("P19").Value = 4
$ws.Range("Q19").Value = 16

# Row 20 (B2)
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 16
$ws.Range("S20").Value = "40.6%"

# Row 21 (B3)
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 18

# Row 22 (B4)
$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 18

# ---------------------------------------------------------------------
# Recorder-list reorderings (same people, new join order) and count bumps
# ---------------------------------------------------------------------
$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G52").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G62").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G74").Value = "Shimaa.ashraf@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G84").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G96").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G106").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G118").Value = "mariam.noureldin@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"

$ws.Range("G120").Value = "basma.hamed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("H120").Value = "140/224"

$ws.Range("G128").Value = "Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G134").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G156").Value = "Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"

# ---------------------------------------------------------------------
# Sessions switching status: Pending -> Recorded (row 24) and
# Pending -> Not Recorded (rows 98, 117, 146, 172), with matching fills
# ---------------------------------------------------------------------

# Row 24: Year2/A2/ANATOMY #1 -> Recorded (copy green style from row 2)
$ws.Range("A2:I2").Copy() | Out-Null
$ws.Range("A24:I24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("G24").Value = "rana.abozaid@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg"
$ws.Range("H24").Value = "92/217"
$ws.Range("I24").Value = "Recorded"

# Row 98: Year2/B1/MICROBIOLOGY #1 -> Not Recorded (copy pink style from row 7)
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A98:I98").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I98").Value = "Not Recorded"

# Row 117: Year2/B2/CARDIOLOGY #1 -> Not Recorded
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A117:I117").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I117").Value = "Not Recorded"

# Row 146: Year2/B3/PHARMACOLOGY #1 -> Not Recorded
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A146:I146").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I146").Value = "Not Recorded"

# Row 172: Year2/B4/PHYSIOLOGY #1 -> Not Recorded
$ws.Range("A7:I7").Copy() | Out-Null
$ws.Range("A172:I172").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("I172").Value = "Not Recorded"
